# CIERRE 31 DIC 2021
# - Amount in the "ARQUITECTO" vale goes from 70000 to 50000 (number + spelled-out text)
# - Active-cell selection on the ARQUITECTO sheet moves from D3 to G8
# - TODAY() driven cells recalc automatically when the workbook is saved

$wb = $excel.ActiveWorkbook

# "ARQUITECTO        " is the tab-selected / active sheet in this workbook
$ws = $wb.ActiveSheet

# Numeric amount
$ws.Range("D1").Value = 50000

# Amount spelled out in words (shared string used by this cell)
$ws.Range("A2").Value = "CINCUENTA      MIL   PESOS 00/100 M.N."

# Move the active cell / selection as captured in the saved view
$ws.Activate()
$ws.Range("G8").Select()
